$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 785005.75
$ws.Range("I125").Value = 385
$ws.Range("J125").Value = 1070322.4
$ws.Range("K125").Value = 3465
$ws.Range("L125").Value = 9632901.6
$ws.Range("M125").Value = -1005
$ws.Range("N125").Value = -9637821.6
$ws.Range("H129").Value = 38462250
$ws.Range("I129").Value = 100000310
$ws.Range("J129").Value = 956.25
$ws.Range("K129").Value = 300000930
$ws.Range("L129").Value = 2868.75
$ws.Range("M129").Value = -299995930
$ws.Range("N129").Value = -12868.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1289.9546
$ws.Range("I2").Value = 1076.2778
$ws.Range("K2").Value = 1076.2778
$ws.Range("M2").Value = -963.2778000000001
$ws.Range("H5").Value = 30100.1
$ws.Range("I5").Value = 33449.832
$ws.Range("J5").Value = 25075.5
$ws.Range("K5").Value = 33449.832
$ws.Range("L5").Value = 25075.5
$ws.Range("M5").Value = -33337.832
$ws.Range("N5").Value = -25299.5
$ws.Range("H8").Value = 5002250
$ws.Range("J8").Value = 4500
$ws.Range("L8").Value = 4500
$ws.Range("N8").Value = -4788
$ws.Range("H11").Value = 4163.3335
$ws.Range("J11").Value = 4163.3335
$ws.Range("L11").Value = 4163.3335
$ws.Range("N11").Value = -4451.3335
$ws.Range("H45").Value = 2121.35
$ws.Range("I45").Value = 983.5625
$ws.Range("J45").Value = 6672.5
$ws.Range("K45").Value = 983.5625
$ws.Range("L45").Value = 6672.5
$ws.Range("M45").Value = -606.5625
$ws.Range("N45").Value = -7426.5
$ws.Range("H102").Value = 1474.875
$ws.Range("I102").Value = 1474.875
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1474.875
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 147.125
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 1289.9546
$ws.Range("I116").Value = 1076.2778
$ws.Range("K116").Value = 1076.2778
$ws.Range("M116").Value = 1217.7222
$ws.Range("H133").Value = 49655.2
$ws.Range("J133").Value = 49655.2
$ws.Range("L133").Value = 49655.2
$ws.Range("N133").Value = -54715.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1289.9546
$ws.Range("I3").Value = 1076.2778
$ws.Range("K3").Value = 1076.2778
$ws.Range("M3").Value = -962.2778000000001
$ws.Range("H4").Value = 30100.1
$ws.Range("I4").Value = 33449.832
$ws.Range("J4").Value = 25075.5
$ws.Range("K4").Value = 33449.832
$ws.Range("L4").Value = 25075.5
$ws.Range("M4").Value = -33334.832
$ws.Range("N4").Value = -25305.5
$ws.Range("H7").Value = 974.75
$ws.Range("I7").Value = 974.75
$ws.Range("K7").Value = 974.75
$ws.Range("M7").Value = -861.75
$ws.Range("H105").Value = 2409.7273
$ws.Range("I105").Value = 2450.7
$ws.Range("K105").Value = 2450.7
$ws.Range("M105").Value = -703.6999999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1706.875
$ws.Range("I132").Value = 1239.5714
$ws.Range("J132").Value = 4978
$ws.Range("K132").Value = 3718.7142
$ws.Range("L132").Value = 14934
$ws.Range("M132").Value = -1188.7142
$ws.Range("N132").Value = -19994
$ws.Range("H135").Value = 37086.156
$ws.Range("J135").Value = 36981
$ws.Range("L135").Value = 36981
$ws.Range("N135").Value = -47121

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I113").Value = 1550.3334
$ws.Range("J113").Value = 766.6667
$ws.Range("K113").Value = 4651.0002
$ws.Range("L113").Value = 2300.0001
$ws.Range("M113").Value = -2481.0002
$ws.Range("N113").Value = -6640.0001
$ws.Range("H131").Value = 1233.3091
$ws.Range("I131").Value = 516.5833
$ws.Range("J131").Value = 1433.3256
$ws.Range("K131").Value = 1549.7499
$ws.Range("L131").Value = 4299.976799999999
$ws.Range("M131").Value = 3490.2501
$ws.Range("N131").Value = -14379.9768

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5000249.5
$ws.Range("I7").Value = 5333333.5
$ws.Range("J7").Value = 4000997.2
$ws.Range("K7").Value = 5333333.5
$ws.Range("L7").Value = 4000997.2
$ws.Range("M7").Value = -5333221.5
$ws.Range("N7").Value = -4001221.2
$ws.Range("H8").Value = 5000249.5
$ws.Range("I8").Value = 5333333.5
$ws.Range("J8").Value = 4000997.2
$ws.Range("K8").Value = 5333333.5
$ws.Range("L8").Value = 4000997.2
$ws.Range("M8").Value = -5333194.5
$ws.Range("N8").Value = -4001275.2
$ws.Range("H11").Value = 5938749.5
$ws.Range("I11").Value = 6638181.5
$ws.Range("J11").Value = 4399999.5
$ws.Range("K11").Value = 6638181.5
$ws.Range("L11").Value = 4399999.5
$ws.Range("M11").Value = -6638042.5
$ws.Range("N11").Value = -4400277.5
$ws.Range("H14").Value = 1018328.8
$ws.Range("I14").Value = 1667221.4
$ws.Range("J14").Value = 44990
$ws.Range("K14").Value = 1667221.4
$ws.Range("L14").Value = 44990
$ws.Range("M14").Value = -1667053.4
$ws.Range("N14").Value = -45326
$ws.Range("H15").Value = 12036.286
$ws.Range("J15").Value = 12036.286
$ws.Range("L15").Value = 12036.286
$ws.Range("N15").Value = -12612.286
$ws.Range("H21").Value = 73111.09
$ws.Range("I21").Value = 78000
$ws.Range("J21").Value = 51111
$ws.Range("K21").Value = 78000
$ws.Range("L21").Value = 51111
$ws.Range("M21").Value = -77827
$ws.Range("N21").Value = -51457
$ws.Range("H24").Value = 76067.5
$ws.Range("J24").Value = 34756.668
$ws.Range("L24").Value = 34756.668
$ws.Range("N24").Value = -35102.668
$ws.Range("H29").Value = 6921.4
$ws.Range("I29").Value = 853.5
$ws.Range("J29").Value = 10966.667
$ws.Range("K29").Value = 853.5
$ws.Range("L29").Value = 10966.667
$ws.Range("M29").Value = -563.5
$ws.Range("N29").Value = -11546.667
$ws.Range("H30").Value = 73111.09
$ws.Range("I30").Value = 78000
$ws.Range("J30").Value = 51111
$ws.Range("K30").Value = 78000
$ws.Range("L30").Value = 51111
$ws.Range("M30").Value = -77895
$ws.Range("N30").Value = -51321
$ws.Range("H36").Value = 2069.8333
$ws.Range("I36").Value = 1200
$ws.Range("J36").Value = 2504.75
$ws.Range("K36").Value = 1200
$ws.Range("L36").Value = 2504.75
$ws.Range("M36").Value = -715
$ws.Range("N36").Value = -3474.75
$ws.Range("H43").Value = 13322.223
$ws.Range("I43").Value = 1450
$ws.Range("J43").Value = 16714.285
$ws.Range("K43").Value = 1450
$ws.Range("L43").Value = 16714.285
$ws.Range("M43").Value = -1299
$ws.Range("N43").Value = -17016.285
$ws.Range("H81").Value = 12036.286
$ws.Range("J81").Value = 12036.286
$ws.Range("L81").Value = 12036.286
$ws.Range("N81").Value = -14032.286
$ws.Range("H84").Value = 12036.286
$ws.Range("J84").Value = 12036.286
$ws.Range("L84").Value = 36108.858
$ws.Range("N84").Value = -46092.858

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 465355.72
$ws.Range("J2").Value = 257490
$ws.Range("L2").Value = 257490
$ws.Range("N2").Value = -257714
$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20452
$ws.Range("H22").Value = 1339.84
$ws.Range("I22").Value = 873.4666999999999
$ws.Range("J22").Value = 2039.4
$ws.Range("K22").Value = 873.4666999999999
$ws.Range("L22").Value = 2039.4
$ws.Range("M22").Value = -578.4666999999999
$ws.Range("N22").Value = -2629.4
$ws.Range("H27").Value = 1339.84
$ws.Range("I27").Value = 873.4666999999999
$ws.Range("J27").Value = 2039.4
$ws.Range("K27").Value = 873.4666999999999
$ws.Range("L27").Value = 2039.4
$ws.Range("M27").Value = -766.4666999999999
$ws.Range("N27").Value = -2253.4
$ws.Range("H46").Value = 26369588
$ws.Range("I46").Value = 62625692
$ws.Range("J46").Value = 1511.5454
$ws.Range("K46").Value = 62625692
$ws.Range("L46").Value = 1511.5454
$ws.Range("M46").Value = -62625504
$ws.Range("N46").Value = -1887.5454
$ws.Range("H55").Value = 212.85715
$ws.Range("I55").Value = 267.77777
$ws.Range("J55").Value = 114
$ws.Range("K55").Value = 267.77777
$ws.Range("L55").Value = 114
$ws.Range("M55").Value = -94.77776999999998
$ws.Range("N55").Value = -460
$ws.Range("H68").Value = 2022.3529
$ws.Range("I68").Value = 1687.5
$ws.Range("J68").Value = 2320
$ws.Range("K68").Value = 1687.5
$ws.Range("L68").Value = 2320
$ws.Range("M68").Value = -938.5
$ws.Range("N68").Value = -3818
$ws.Range("H71").Value = 2022.3529
$ws.Range("I71").Value = 1687.5
$ws.Range("J71").Value = 2320
$ws.Range("K71").Value = 8437.5
$ws.Range("L71").Value = 11600
$ws.Range("M71").Value = -4693.5
$ws.Range("N71").Value = -19088

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 781.8182
$ws.Range("J14").Value = 2333.3333
$ws.Range("L14").Value = 2333.3333
$ws.Range("N14").Value = -2669.3333
